$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 24.35712233333334
$ws.Range("H2").Value = 73.07136700000001
$ws.Range("I2").Value = 0.3750500562097488
$ws.Range("J2").Value = 0.3750500562097488
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 3548.850861903765
$ws.Range("R2").Value = 31939.65775713389
$ws.Range("S2").Value = 0.1074870114270837
$ws.Range("T2").Value = 0.1074870114270837
$ws.Range("G3").Value = 24.35712233333334
$ws.Range("H3").Value = 73.07136700000001
$ws.Range("I3").Value = 0.3750500562097488
$ws.Range("J3").Value = 0.3750500562097488
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 4111.475007682294
$ws.Range("R3").Value = 37003.27506914065
$ws.Range("S3").Value = 0.1245276790515352
$ws.Range("T3").Value = 0.1245276790515351
$ws.Range("G4").Value = 24.35712233333334
$ws.Range("H4").Value = 73.07136700000001
$ws.Range("I4").Value = 0.3750500562097488
$ws.Range("J4").Value = 0.3750500562097488
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 3120.784423315588
$ws.Range("R4").Value = 28087.05980984029
$ws.Range("S4").Value = 0.09452180551493795
$ws.Range("T4").Value = 0.09452180551493793
$ws.Range("G5").Value = 24.35712233333334
$ws.Range("H5").Value = 73.07136700000001
$ws.Range("I5").Value = 0.3750500562097488
$ws.Range("J5").Value = 0.3750500562097488
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1601.750645974998
$ws.Range("R5").Value = 14415.75581377498
$ws.Range("S5").Value = 0.04851356021619205
$ws.Range("T5").Value = 0.04851356021619204
$ws.Range("I6").Value = 0.2805618708302703
$ws.Range("J6").Value = 0.2805618708302702
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 2654.771598158359
$ws.Range("R6").Value = 23892.94438342523
$ws.Range("S6").Value = 0.08040728568527908
$ws.Range("T6").Value = 0.08040728568527906
$ws.Range("I7").Value = 0.2805618708302703
$ws.Range("J7").Value = 0.2805618708302702
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.09315481500770938
$ws.Range("T7").Value = 0.09315481500770936
$ws.Range("I8").Value = 0.2805618708302703
$ws.Range("J8").Value = 0.2805618708302702
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 2334.550020101099
$ws.Range("R8").Value = 21010.95018090989
$ws.Range("S8").Value = 0.07070846717776501
$ws.Range("T8").Value = 0.070708467177765
$ws.Range("I9").Value = 0.2805618708302703
$ws.Range("J9").Value = 0.2805618708302702
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 1198.213812790406
$ws.Range("R9").Value = 10783.92431511365
$ws.Range("S9").Value = 0.03629130295951683
$ws.Range("T9").Value = 0.03629130295951682
$ws.Range("G10").Value = 22.31748066666667
$ws.Range("H10").Value = 66.952442
$ws.Range("I10").Value = 0.3436437303202491
$ws.Range("J10").Value = 0.343643730320249
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 3251.67355221727
$ws.Range("R10").Value = 29265.06196995543
$ws.Range("S10").Value = 0.09848615393119936
$ws.Range("T10").Value = 0.09848615393119935
$ws.Range("G11").Value = 22.31748066666667
$ws.Range("H11").Value = 66.952442
$ws.Range("I11").Value = 0.3436437303202491
$ws.Range("J11").Value = 0.343643730320249
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 3767.184100802416
$ws.Range("R11").Value = 33904.65690722175
$ws.Range("S11").Value = 0.1140998526699592
$ws.Range("T11").Value = 0.1140998526699592
$ws.Range("G12").Value = 22.31748066666667
$ws.Range("H12").Value = 66.952442
$ws.Range("I12").Value = 0.3436437303202491
$ws.Range("J12").Value = 0.343643730320249
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 2859.452979667677
$ws.Range("R12").Value = 25735.07681700909
$ws.Range("S12").Value = 0.0866066417160933
$ws.Range("T12").Value = 0.08660664171609328
$ws.Range("G13").Value = 22.31748066666667
$ws.Range("H13").Value = 66.952442
$ws.Range("I13").Value = 0.3436437303202491
$ws.Range("J13").Value = 0.343643730320249
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 1467.621609201639
$ws.Range("R13").Value = 13208.59448281475
$ws.Range("S13").Value = 0.04445108200299722
$ws.Range("T13").Value = 0.04445108200299722
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04834033333333334
$ws.Range("H14").Value = 0.145021
$ws.Range("I14").Value = 0.0007443426397318391
$ws.Range("J14").Value = 0.0007443426397318388
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 7.043222564101557
$ws.Range("R14").Value = 63.38900307691402
$ws.Range("S14").Value = 0.0002133239670220911
$ws.Range("T14").Value = 0.000213323967022091
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04834033333333334
$ws.Range("H15").Value = 0.145021
$ws.Range("I15").Value = 0.0007443426397318391
$ws.Range("J15").Value = 0.0007443426397318388
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 8.159833893474225
$ws.Range("R15").Value = 73.43850504126802
$ws.Range("S15").Value = 0.0002471437073803844
$ws.Range("T15").Value = 0.0002471437073803844
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04834033333333334
$ws.Range("H16").Value = 0.145021
$ws.Range("I16").Value = 0.0007443426397318391
$ws.Range("J16").Value = 0.0007443426397318388
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 6.19366102530489
$ws.Range("R16").Value = 55.74294922774401
$ws.Range("S16").Value = 0.000187592586814228
$ws.Range("T16").Value = 0.000187592586814228
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.04834033333333334
$ws.Range("H17").Value = 0.145021
$ws.Range("I17").Value = 0.0007443426397318391
$ws.Range("J17").Value = 0.0007443426397318388
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 3.178912479219667
$ws.Range("R17").Value = 28.610212312977
$ws.Range("S17").Value = 0.00009628237851513559
$ws.Range("T17").Value = 0.00009628237851513556

Write-Host "Updated 182 cells"
